$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'" + '61.445.97'
$ws.Cells.Item(2, 5).Value = '  +0.86%  '
$ws.Cells.Item(3, 4).Value = "'" + '3.436.06'
$ws.Cells.Item(3, 5).Value = '  +1.53%  '
$ws.Cells.Item(5, 4).Value = "'" + '575.44'
$ws.Cells.Item(5, 5).Value = '  +0.61%  '
$ws.Cells.Item(6, 4).Value = "'" + '145.48'
$ws.Cells.Item(6, 5).Value = '  +6.36%  '
$ws.Cells.Item(7, 4).Value = "'" + '3.437.21'
$ws.Cells.Item(7, 5).Value = '  +1.61%  '
$ws.Cells.Item(9, 4).Value = "'" + '0.477'
$ws.Cells.Item(9, 5).Value = '  +1.81%  '
$ws.Cells.Item(10, 4).Value = "'" + '7.67'
$ws.Cells.Item(10, 5).Value = '  +0.59%  '
$ws.Cells.Item(11, 5).Value = '  +3.30%  '
$ws.Cells.Item(12, 5).Value = '  +1.63%  '
$ws.Cells.Item(13, 4).Value = "'" + '4.022.09'
$ws.Cells.Item(13, 5).Value = '  +1.56%  '
$ws.Cells.Item(14, 4).Value = "'" + '28.02'
$ws.Cells.Item(14, 5).Value = '  +6.15%  '
$ws.Cells.Item(15, 5).Value = '  -0.55%  '
$ws.Cells.Item(16, 5).Value = '  +1.53%  '
$ws.Cells.Item(17, 4).Value = "'" + '3.431.26'
$ws.Cells.Item(17, 5).Value = '  +1.49%  '
$ws.Cells.Item(18, 4).Value = "'" + '61.544.61'
$ws.Cells.Item(18, 5).Value = '  +0.94%  '
$ws.Cells.Item(19, 4).Value = "'" + '6.28'
$ws.Cells.Item(19, 5).Value = '  +7.76%  '
$ws.Cells.Item(20, 4).Value = "'" + '14.24'
$ws.Cells.Item(20, 5).Value = '  +2.88%  '
$ws.Cells.Item(21, 4).Value = "'" + '9.41'
$ws.Cells.Item(21, 5).Value = '  +1.35%  '
$ws.Cells.Item(22, 4).Value = "'" + '395.65'
$ws.Cells.Item(22, 5).Value = '  +5.62%  '
$ws.Cells.Item(23, 4).Value = "'" + '0.566'
$ws.Cells.Item(23, 5).Value = '  +3.09%  '
$ws.Cells.Item(24, 4).Value = "'" + '73.79'
$ws.Cells.Item(24, 5).Value = '  +4.37%  '
$ws.Cells.Item(25, 4).Value = "'" + '0.997'
$ws.Cells.Item(25, 5).Value = '  -0.34%  '
$ws.Cells.Item(26, 4).Value = "'" + '5.72'
$ws.Cells.Item(26, 5).Value = '  -0.11%  '
$ws.Cells.Item(27, 4).Value = "'" + '0.0000124'
$ws.Cells.Item(27, 5).Value = '  +1.01%  '
$ws.Cells.Item(28, 4).Value = "'" + '3.574.14'
$ws.Cells.Item(28, 5).Value = '  +1.86%  '
$ws.Cells.Item(29, 5).Value = '  +4.63%  '
$ws.Cells.Item(30, 4).Value = "'" + '7.60'
$ws.Cells.Item(30, 5).Value = '  +3.42%  '
$ws.Cells.Item(31, 5).Value = '  +0.17%  '
$ws.Cells.Item(32, 4).Value = "'" + '8.24'
$ws.Cells.Item(32, 5).Value = '  +2.28%  '
$ws.Cells.Item(33, 4).Value = "'" + '1.47'
$ws.Cells.Item(33, 5).Value = '  -9.00%  '
$ws.Cells.Item(34, 5).Value = '  +2.02%  '
$ws.Cells.Item(35, 5).Value = '  -0.04%  '
$ws.Cells.Item(36, 4).Value = "'" + '23.91'
$ws.Cells.Item(36, 5).Value = '  +2.42%  '
$ws.Cells.Item(37, 2).Value = 'RenzoRestakedETH'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Cells.Item(37, 4).Value = "'" + '3.464.13'
$ws.Cells.Item(37, 5).Value = '  +1.84%  '
$ws.Cells.Item(38, 2).Value = 'Aptos'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(38, 4).Value = "'" + '7.01'
$ws.Cells.Item(38, 5).Value = '  +3.17%  '
$ws.Cells.Item(39, 2).Value = 'NEARProtocol'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(39, 4).Value = "'" + '5.11'
$ws.Cells.Item(39, 5).Value = '  -0.01%  '
$ws.Cells.Item(40, 2).Value = 'ImmutableX'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(40, 4).Value = "'" + '1.55'
$ws.Cells.Item(40, 5).Value = '  +0.59%  '
$ws.Cells.Item(41, 4).Value = "'" + '167.61'
$ws.Cells.Item(41, 5).Value = '  +1.65%  '
$ws.Cells.Item(42, 4).Value = "'" + '0.0784'
$ws.Cells.Item(42, 5).Value = '  +2.45%  '
$ws.Cells.Item(43, 4).Value = "'" + '26.98'
$ws.Cells.Item(43, 5).Value = '  +4.97%  '
$ws.Cells.Item(44, 5).Value = '  +3.46%  '
$ws.Cells.Item(45, 2).Value = 'Stacks'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(45, 4).Value = "'" + '1.74'
$ws.Cells.Item(45, 5).Value = '  -0.28%  '
$ws.Cells.Item(46, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(46, 4).Value = "'" + '1.00'
$ws.Cells.Item(46, 5).Value = '  +0.05%  '
$ws.Cells.Item(47, 5).Value = '  +3.16%  '
$ws.Cells.Item(48, 4).Value = "'" + '42.31'
$ws.Cells.Item(48, 5).Value = '  +0.93%  '
$ws.Cells.Item(49, 4).Value = "'" + '2.601.31'
$ws.Cells.Item(49, 5).Value = '  +3.47%  '
$ws.Cells.Item(50, 5).Value = '  -1.31%  '
$ws.Cells.Item(51, 5).Value = '  +2.43%  '
